# Apply "Fixed update to excel issue" edit:
# The forecast table rolls forward by one week: each row's Week_Start_Date
# and forecast values (MyForecast / Amazon Mean / P70 / P80 / P90) are
# updated to reflect the refreshed forecast, and the Summary sheet's
# derived metrics are refreshed to match.

$wb = $excel.ActiveWorkbook
$fc = $wb.Worksheets.Item("Forecast Comparison")
$sm = $wb.Worksheets.Item("Summary")

$fc.Range("B2").Value = "'2025-02-02"
$fc.Range("D2").Value = 35
$fc.Range("E2").Value = 31
$fc.Range("F2").Value = 37
$fc.Range("G2").Value = 42
$fc.Range("H2").Value = 50
$fc.Range("B3").Value = "'2025-02-09"
$fc.Range("D3").Value = 34
$fc.Range("E3").Value = 32
$fc.Range("F3").Value = 39
$fc.Range("G3").Value = 45
$fc.Range("H3").Value = 54
$fc.Range("B4").Value = "'2025-02-16"
$fc.Range("D4").Value = 39
$fc.Range("E4").Value = 33
$fc.Range("F4").Value = 40
$fc.Range("G4").Value = 46
$fc.Range("H4").Value = 56
$fc.Range("B5").Value = "'2025-02-23"
$fc.Range("D5").Value = 43
$fc.Range("E5").Value = 33
$fc.Range("F5").Value = 39
$fc.Range("G5").Value = 46
$fc.Range("H5").Value = 57
$fc.Range("B6").Value = "'2025-03-02"
$fc.Range("D6").Value = 43
$fc.Range("E6").Value = 33
$fc.Range("F6").Value = 40
$fc.Range("G6").Value = 47
$fc.Range("H6").Value = 58
$fc.Range("B7").Value = "'2025-03-09"
$fc.Range("D7").Value = 43
$fc.Range("E7").Value = 33
$fc.Range("F7").Value = 40
$fc.Range("G7").Value = 47
$fc.Range("H7").Value = 58
$fc.Range("B8").Value = "'2025-03-16"
$fc.Range("D8").Value = 40
$fc.Range("E8").Value = 32
$fc.Range("F8").Value = 39
$fc.Range("G8").Value = 46
$fc.Range("H8").Value = 58
$fc.Range("B9").Value = "'2025-03-23"
$fc.Range("D9").Value = 38
$fc.Range("E9").Value = 31
$fc.Range("F9").Value = 38
$fc.Range("G9").Value = 45
$fc.Range("H9").Value = 57
$fc.Range("B10").Value = "'2025-03-30"
$fc.Range("D10").Value = 39
$fc.Range("E10").Value = 30
$fc.Range("F10").Value = 37
$fc.Range("G10").Value = 45
$fc.Range("H10").Value = 56
$fc.Range("B11").Value = "'2025-04-06"
$fc.Range("D11").Value = 40
$fc.Range("E11").Value = 31
$fc.Range("F11").Value = 38
$fc.Range("G11").Value = 46
$fc.Range("H11").Value = 59
$fc.Range("B12").Value = "'2025-04-13"
$fc.Range("D12").Value = 39
$fc.Range("E12").Value = 30
$fc.Range("F12").Value = 36
$fc.Range("G12").Value = 44
$fc.Range("H12").Value = 56
$fc.Range("B13").Value = "'2025-04-20"
$fc.Range("D13").Value = 39
$fc.Range("E13").Value = 30
$fc.Range("F13").Value = 36
$fc.Range("G13").Value = 44
$fc.Range("H13").Value = 55
$fc.Range("B14").Value = "'2025-04-27"
$fc.Range("D14").Value = 38
$fc.Range("E14").Value = 29
$fc.Range("F14").Value = 36
$fc.Range("G14").Value = 43
$fc.Range("H14").Value = 55
$fc.Range("B15").Value = "'2025-05-04"
$fc.Range("D15").Value = 36
$fc.Range("E15").Value = 28
$fc.Range("F15").Value = 34
$fc.Range("G15").Value = 42
$fc.Range("H15").Value = 54
$fc.Range("B16").Value = "'2025-05-11"
$fc.Range("D16").Value = 36
$fc.Range("E16").Value = 28
$fc.Range("F16").Value = 34
$fc.Range("G16").Value = 41
$fc.Range("H16").Value = 53
$fc.Range("B17").Value = "'2025-05-18"
$fc.Range("D17").Value = 36
$fc.Range("E17").Value = 28
$fc.Range("F17").Value = 34
$fc.Range("G17").Value = 42
$fc.Range("H17").Value = 55

# Summary sheet updates
$sm.Range("B2").Value = "2022-12-25 to 2025-01-26"
$sm.Range("B8").Value = "8925 units"
$sm.Range("B9").Value = "'619"
$sm.Range("B10").Value = "'315"
$sm.Range("B11").Value = "'151"
$sm.Range("B12").Value = "'43"
$sm.Range("B14").Value = "'34"
